$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '(0   18.0  26.0  34.0  44.0)'
$ws.Range("E2").Value = '(0   18.0  26.0  34.0  44.0)'
$ws.Range("D3").Value = '(0   10.0  18.0  26.0   0.0)'
$ws.Range("E3").Value = '(0   10.0  18.0  26.0   0.0)'
$ws.Range("D11").Value = '(0  1  2  3  4  5)'
$ws.Range("E11").Value = '(0  1  2  3  4  5)'
$ws.Range("D12").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("E12").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("D13").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("E13").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("D14").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("E14").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("D15").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("E15").Value = '(0.2  0.2  0.2  0.2  0.2  0.2)'
$ws.Range("D16").Value = '(1.9  1.9  1.9  1.9  1.9  1.9)'
$ws.Range("E16").Value = '(1.9  1.9  1.9  1.9  1.9  1.9)'
$ws.Range("D17").Value = '(0  1  2  3  4  5)'
$ws.Range("E17").Value = '(0  1  2  3  4  5)'
$ws.Range("D18").Value = '(2  2  2  2  2  2)'
$ws.Range("E18").Value = '(2  2  2  2  2  2)'
$ws.Range("D19").Value = '(2  2  2  2  2  2)'
$ws.Range("E19").Value = '(2  2  2  2  2  2)'
$ws.Range("D20").Value = '(0  0.5  0.8 1.0)'
$ws.Range("E20").Value = '(0  0.5  0.8 1.0)'
$ws.Range("D21").Value = '(0  0.5  1.0 1.0)'
$ws.Range("E21").Value = '(0  0.5  1.0 1.0)'
$ws.Range("D22").Value = '(0  0.5  0.8 1.0)'
$ws.Range("E22").Value = '(0  0.5  0.8 1.0)'
$ws.Range("D23").Value = '(0  0.5  1.0 1.0)'
$ws.Range("E23").Value = '(0  0.5  1.0 1.0)'
$ws.Range("C24").Value = 'Fraction of N range in plant (e.g.  0.5 is 50% between N min and N max)'
$ws.Range("D24").Value = '(0  0.5  0.8 1.0)'
$ws.Range("E24").Value = '(0  0.5  0.8 1.0)'
$ws.Range("D25").Value = '(0  0.5  1.0 1.0)'
$ws.Range("E25").Value = '(0  0.5  1.0 1.0)'
$ws.Range("C26").Value = 'Fraction of N range in plant (e.g.  0.5 is 50% between N min and N max)'
$ws.Range("D26").Value = '(0  0.5  0.8 1.0)'
$ws.Range("E26").Value = '(0  0.5  0.8 1.0)'
$ws.Range("D27").Value = '(0  0.5  1.0 1.0)'
$ws.Range("E27").Value = '(0  0.5  1.0 1.0)'
$ws.Range("D30").Value = '(0  5  10  15  20)'
$ws.Range("E30").Value = '(0  5  10  15  20)'
$ws.Range("D31").Value = '(0  0.2  0.5  1.0)'
$ws.Range("E31").Value = '(0  0.2  0.5  1.0)'
$ws.Range("D32").Value = '(0  280  350  550  1200)'
$ws.Range("E32").Value = '(0  280  350  550  1200)'
$ws.Range("D33").Value = '(0  0.9  1  1.2  1.5)'
$ws.Range("E33").Value = '(0  0.9  1  1.2  1.5)'
$ws.Range("D34").Value = '(0  0.9  1  1.2  1.5)'
$ws.Range("E34").Value = '(0  0.9  1  1.2  1.5)'
$ws.Cells.Select()
